# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gains a new "property_category" column
# (with value "stock" for every data row), inserted right before the
# existing "date" column. Everything that used to live in/after that
# column (date, legislator_name, legislator_id) shifts one column to
# the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (pushing date/legislator_name/legislator_id,
# formerly H/I/J, to I/J/K) while preserving formatting on both sides.
$ws.Range("H1:H5").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "property_category"

# Populate the new column for each of the 4 stock rows.
$ws.Range("H2:H5").Value = "stock"
